# today my new commit
# Add a new "ApplelaptopProduct" row to the Customer sheet and update the
# selection/window view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer")

# Add the new data row (row 5): A5 = "Campaign Name", B5 = "ApplelaptopProduct"
$ws.Range("A5").Value = "Campaign Name"
$ws.Range("B5").Value = "ApplelaptopProduct"

# Move the active selection to the newly added cell
$ws.Range("B5").Select()

# Adjust the workbook window height (view state)
$excel.ActiveWindow.Height = 6180
